$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text values (also rewrites shared strings table) ---
$ws.Range("A1").Value = "WCenter"
$ws.Range("B1").Value = "WRType"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Eng"
$ws.Range("E1").Value = "WorkNum"
$ws.Range("F1").Value = "CreationTime"
$ws.Range("G1").Value = "Description"

# --- Clear old sample values from row 2 ---
$ws.Range("B2:E2").ClearContents()

# --- Header styling: bold 14pt black Arial, green fill, thin borders, centered ---
$hdr = $ws.Range("A1:G1")
$hdr.Font.Bold = $true
$hdr.Font.Size = 14
$hdr.Font.Name = "Arial"
$hdr.Font.Color = 0
$hdr.Interior.Color = 13434828
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

# --- Row height / column widths ---
$ws.Rows.Item(1).RowHeight = 40
$ws.Columns.Item(1).ColumnWidth = 23.285714285714285
$ws.Columns.Item(2).ColumnWidth = 17.285714285714285
$ws.Columns.Item(3).ColumnWidth = 17.285714285714285
$ws.Columns.Item(4).ColumnWidth = 17.285714285714285
$ws.Columns.Item(5).ColumnWidth = 17.285714285714285
$ws.Columns.Item(6).ColumnWidth = 17.285714285714285
$ws.Columns.Item(7).ColumnWidth = 29.285714285714285

# --- View: hide gridlines, 50% zoom ---
$excel.ActiveWindow.DisplayGridlines = $false
$excel.ActiveWindow.Zoom = 50

# --- Page setup: landscape ---
$ws.PageSetup.Orientation = 2
